# Auto-generated edit script applying numeric corrections to Sheets per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 224
$ws.Range("I5").Value = 115.55556
$ws.Range("K5").Value = 115.55556
$ws.Range("M5").Value = -0.5555599999999998
$ws.Range("H18").Value = 711.7778
$ws.Range("I18").Value = 711.7778
$ws.Range("K18").Value = 711.7778
$ws.Range("M18").Value = -427.7778
$ws.Range("H40").Value = 5077.9
$ws.Range("J40").Value = 6957.2
$ws.Range("L40").Value = 6957.2
$ws.Range("N40").Value = -7307.2
$ws.Range("H94").Value = 6667.3335
$ws.Range("I94").Value = 6667.3335
$ws.Range("K94").Value = 6667.3335
$ws.Range("M94").Value = -6216.3335
$ws.Range("H100").Value = 1526.9062
$ws.Range("I100").Value = 1091.08
$ws.Range("J100").Value = 3083.4285
$ws.Range("K100").Value = 1091.08
$ws.Range("L100").Value = 3083.4285
$ws.Range("M100").Value = -550.0799999999999
$ws.Range("N100").Value = -4165.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7793.7856
$ws.Range("I32").Value = 4838.2856
$ws.Range("K32").Value = 4838.2856
$ws.Range("M32").Value = -4551.2856
$ws.Range("H122").Value = 4841.857
$ws.Range("I122").Value = 4844.4546
$ws.Range("J122").Value = 4832.3335
$ws.Range("K122").Value = 14533.3638
$ws.Range("L122").Value = 14497.0005
$ws.Range("M122").Value = -12083.3638
$ws.Range("N122").Value = -19397.0005
$ws.Range("H132").Value = 25541.605
$ws.Range("I132").Value = 27156.975
$ws.Range("K132").Value = 81470.92499999999
$ws.Range("M132").Value = -78940.92499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 251174.25
$ws.Range("I22").Value = 500899.5
$ws.Range("K22").Value = 500899.5
$ws.Range("M22").Value = -500726.5
$ws.Range("H86").Value = 2192.4546
$ws.Range("I86").Value = 2056.4285
$ws.Range("J86").Value = 2430.5
$ws.Range("K86").Value = 2056.4285
$ws.Range("L86").Value = 2430.5
$ws.Range("M86").Value = -933.4285
$ws.Range("N86").Value = -4676.5
$ws.Range("H89").Value = 2192.4546
$ws.Range("I89").Value = 2056.4285
$ws.Range("J89").Value = 2430.5
$ws.Range("K89").Value = 10282.1425
$ws.Range("L89").Value = 12152.5
$ws.Range("M89").Value = -4666.1425
$ws.Range("N89").Value = -23384.5
$ws.Range("H105").Value = 4804.8335
$ws.Range("I105").Value = 4678.2144
$ws.Range("K105").Value = 4678.2144
$ws.Range("M105").Value = -2931.2144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 494.63635
$ws.Range("I22").Value = 303.75
$ws.Range("J22").Value = 723.7
$ws.Range("K22").Value = 303.75
$ws.Range("L22").Value = 723.7
$ws.Range("M22").Value = 46.25
$ws.Range("N22").Value = -1423.7
$ws.Range("H62").Value = 4124.25
$ws.Range("I62").Value = 4165.8335
$ws.Range("J62").Value = 3999.5
$ws.Range("K62").Value = 4165.8335
$ws.Range("L62").Value = 3999.5
$ws.Range("M62").Value = -3541.8335
$ws.Range("N62").Value = -5247.5
$ws.Range("H65").Value = 4124.25
$ws.Range("I65").Value = 4165.8335
$ws.Range("J65").Value = 3999.5
$ws.Range("K65").Value = 20829.1675
$ws.Range("L65").Value = 19997.5
$ws.Range("M65").Value = -17709.1675
$ws.Range("N65").Value = -26237.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2990.48
$ws.Range("I102").Value = 2523.35
$ws.Range("K102").Value = 2523.35
$ws.Range("M102").Value = -901.3499999999999
$ws.Range("H107").Value = 72675.21000000001
$ws.Range("I107").Value = 166959.17
$ws.Range("J107").Value = 1962.25
$ws.Range("K107").Value = 166959.17
$ws.Range("L107").Value = 1962.25
$ws.Range("M107").Value = -165039.17
$ws.Range("N107").Value = -5802.25
$ws.Range("H113").Value = 1166.8462
$ws.Range("I113").Value = 1166.8462
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1166.8462
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1003.1538
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 4640.5454
$ws.Range("I126").Value = 3881.7727
$ws.Range("J126").Value = 6158.091
$ws.Range("K126").Value = 11645.3181
$ws.Range("L126").Value = 18474.273
$ws.Range("M126").Value = -9175.3181
$ws.Range("N126").Value = -23414.273
$ws.Range("H132").Value = 28540.264
$ws.Range("I132").Value = 33882.355
$ws.Range("K132").Value = 101647.065
$ws.Range("M132").Value = -99117.065

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 14471.37
$ws.Range("J46").Value = 5131.6665
$ws.Range("L46").Value = 5131.6665
$ws.Range("N46").Value = -5507.6665
$ws.Range("H55").Value = 1356.6562
$ws.Range("I55").Value = 1162.76
$ws.Range("J55").Value = 2049.1428
$ws.Range("K55").Value = 1162.76
$ws.Range("L55").Value = 2049.1428
$ws.Range("M55").Value = -989.76
$ws.Range("N55").Value = -2395.1428
$ws.Range("H68").Value = 3084.2
$ws.Range("I68").Value = 1944
$ws.Range("J68").Value = 4794.5
$ws.Range("K68").Value = 1944
$ws.Range("L68").Value = 4794.5
$ws.Range("M68").Value = -1195
$ws.Range("N68").Value = -6292.5
$ws.Range("H71").Value = 3084.2
$ws.Range("I71").Value = 1944
$ws.Range("J71").Value = 4794.5
$ws.Range("K71").Value = 9720
$ws.Range("L71").Value = 23972.5
$ws.Range("M71").Value = -5976
$ws.Range("N71").Value = -31460.5
$ws.Range("H93").Value = 1165.5
$ws.Range("I93").Value = 832
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 832
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = 416
$ws.Range("N93").Value = -5996
$ws.Range("H132").Value = 35059
$ws.Range("I132").Value = 41038.16
$ws.Range("J132").Value = 4166.6665
$ws.Range("K132").Value = 123114.48
$ws.Range("L132").Value = 12499.9995
$ws.Range("M132").Value = -120584.48
$ws.Range("N132").Value = -17559.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 46186.824
$ws.Range("I132").Value = 46186.824
$ws.Range("K132").Value = 138560.472
$ws.Range("M132").Value = -136030.472
